# merge legacy sni spec to sni page
# - remove RTOS term
# - rename virtual machine to core engine
# - remove MicroEJ prefix
# - sync startup code snippet from latest microej_main cco

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "The Java [RTOS task]" box -> "Core Engine [Task]"
$shape1 = $s.Shapes.Item(14)
$tr1 = $shape1.TextFrame.TextRange
[void]$tr1.Replace("The Java ", "Core Engine ")
[void]$tr1.Replace("RTOS task", "Task")

# "Another C [RTOS task]" box -> "Another C [Task]"
$shape2 = $s.Shapes.Item(15)
$tr2 = $shape2.TextFrame.TextRange
[void]$tr2.Replace("RTOS task", "Task")
